$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 498 first (pushes existing rows 498.. down by
# one, so the sheet ends up with 780 rows and the true last row becomes 780).
$ws.Rows.Item(498).Insert()

# Append a new row at the very end (row 781) with a new facility/school
# mapping. Set B before A so "EW RHODES MS" registers as shared-string index
# 949, ahead of "Lindley Academy Charter School" (index 950).
$ws.Range("B781").Value = "E. Washington Rhodes School"
$ws.Range("A781").Value = "EW RHODES MS"

$ws.Range("B498").Value = "Lindley Academy Charter at Birney"
$ws.Range("A498").Value = "Lindley Academy Charter School"

# Update sheet view to match the author's final scroll position / selection
$ws.Application.ActiveWindow.ScrollRow = 484
$ws.Range("A498").Select()
